$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.984.70"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "3.508.24"
$ws.Range("E3").Value = "  -2.17%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "202.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "550.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.59%  "
$ws.Range("D7").Value = "3.501.06"
$ws.Range("E7").Value = "  -2.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.653"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +11.60%  "
$ws.Range("E12").Value = "  -4.85%  "
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").Value = "4.084.03"
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("D16").Value = "3.505.55"
$ws.Range("E16").Value = "  -2.19%  "
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").Value = "66.793.17"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.02%  "
$ws.Range("E21").Value = "  -3.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -10.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.09%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "687.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "63.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("E36").Value = "  -4.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "39.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.404"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.42%  "
$ws.Range("D41").Value = "3.117.99"
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("E43").Value = "  -3.76%  "
$ws.Range("D44").Value = "0.0₃0707"
$ws.Range("E44").Value = "  -10.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +17.75%  "
$ws.Range("E46").Value = "  -12.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0396"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.127"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.75%  "
